# Weekly update: insert a new data row (new date) for
# "Hortaliza, Terminal La Palmera de La Serena - Cilantro" ahead of the
# existing row 111, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 111; everything currently at
# row 111 and below (111..160) shifts down to 112..161.
$ws.Rows.Item(111).Insert()

# Populate the newly inserted row 111 with the new weekly observation.
$ws.Cells.Item(111, 1).Value2 = 8
$ws.Cells.Item(111, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(111, 3).Value2 = "Coquimbo"
$ws.Cells.Item(111, 4).Value2 = 44784
$ws.Cells.Item(111, 5).Value2 = 4
$ws.Cells.Item(111, 6).Value2 = 100112040
$ws.Cells.Item(111, 7).Value2 = "Cilantro"
$ws.Cells.Item(111, 8).Value2 = "Sin especificar"
$ws.Cells.Item(111, 9).Value2 = "Primera"
$ws.Cells.Item(111, 10).Value2 = 2600
$ws.Cells.Item(111, 11).Value2 = 2000
$ws.Cells.Item(111, 12).Value2 = 2500
$ws.Cells.Item(111, 13).Value2 = 2250
$ws.Cells.Item(111, 14).Value2 = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(111, 15).Value2 = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(111, 16).Value2 = 1500
$ws.Cells.Item(111, 17).Value2 = 1.5
$ws.Cells.Item(111, 18).Value2 = "Hortaliza"
